# Updated symbol list on Mon Dec 19 09:58:13 UTC 2022 with GitHub Actions
#
# Applies the per-cell text updates from the scraped diff to Sheet1 of the
# crypto-ranking workbook: refreshed "Price" (column D) quotes for most
# rows, a couple of "Volume(1h)" (column E) label tweaks, and two pairs of
# rows (14/15 and 42/43) whose Coin/Link/Price/Volume values were
# re-ordered between runs.
#
# Column D stores prices as literal text (t="inlineStr") so values like
# "247.50" or "0.04700" keep their exact trailing zeros instead of being
# normalised as numbers. Plain `.Value = "<numeric-looking string>"` would
# let Excel auto-convert that into a real number and drop the formatting,
# so for those cells we first mark the cell as Text ("@") before writing
# the string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Column D "Price" refreshes (plain row updates) ---
Set-TextValue "D2"  "247.50"
Set-TextValue "D3"  "21.73"
Set-TextValue "D4"  "5.478"
Set-TextValue "D5"  "0.05694"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "0.8064"
Set-TextValue "D8"  "1.038"
Set-TextValue "D9"  "0.1476"
Set-TextValue "D10" "0.07338"
Set-TextValue "D11" "0.03153"
Set-TextValue "D12" "0.02957"
Set-TextValue "D13" "0.09293"

# --- Rows 14 & 15 swap places (BitForexToken <-> MCDex) ---
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "3.446"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001654"
$ws.Range("E15").Value = "14BitForexTokenBF"

# --- More column D "Price" refreshes ---
Set-TextValue "D16" "0.04702"
Set-TextValue "D17" "0.0005870"
Set-TextValue "D18" "0.006350"

Set-TextValue "D19" "0.005048"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"

Set-TextValue "D20" "0.001047"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "0.0003135"
Set-TextValue "D23" "3.775"
Set-TextValue "D24" "6.433"
Set-TextValue "D25" "2.111"

Set-TextValue "D40" "0.04101"

Set-TextValue "D41" "0.006935"
$ws.Range("E41").Value = "40KickTokenKICK"

# --- Rows 42 & 43 swap places (CEJI <-> BKEXToken) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1045"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002972"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining column D "Price" refreshes ---
Set-TextValue "D44" "0.008136"
Set-TextValue "D45" "0.00005835"
Set-TextValue "D47" "0.0005500"
Set-TextValue "D48" "0.6825"
Set-TextValue "D49" "0.009481"
